# "added 4wk low sales check"
# Recomputes the MyForecast / Inventory Coverage / Stockout Risk /
# Reorder Urgency / Seasonality Index columns on "Forecast Comparison"
# (rows 2-17) and rolls the new totals up into "Summary" (B9:B12, B14).

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row => MyForecast(D), Inventory Coverage(H), Stockout Risk(I),
#        Reorder Urgency(J), Seasonality Index(L)
$rows = @(
    @{ Row = 2;  D = 2; H = 2.69; I = "Low";  J = "Normal"; L = 1.03 },
    @{ Row = 3;  D = 1; H = 2.75; I = "Low";  J = "Normal"; L = 1.07 },
    @{ Row = 4;  D = 0; H = 4.67; I = "Low";  J = "Normal"; L = 1.11 },
    @{ Row = 5;  D = 0; H = 5.5;  I = "Low";  J = "Normal"; L = 0.87 },
    @{ Row = 6;  D = 0; H = 4.5;  I = "Low";  J = "Normal"; L = 1.08 },
    @{ Row = 7;  D = 0; H = 3.5;  I = "Low";  J = "Normal"; L = 1.19 },
    @{ Row = 8;  D = 0; H = 2.5;  I = "Low";  J = "Normal"; L = 1.2  },
    @{ Row = 9;  D = 0; H = 1.5;  I = "Low";  J = "Normal"; L = 1.08 },
    @{ Row = 10; D = 0; H = 0.5;  I = "Low";  J = "Urgent"; L = 1.09 },
    @{ Row = 11; D = 0; H = 0;    I = "High"; J = "Urgent"; L = 0.86 },
    @{ Row = 12; D = 0; H = 0;    I = "High"; J = "Urgent"; L = 1.08 },
    @{ Row = 13; D = 0; H = 0;    I = "High"; J = "Urgent"; L = 0.98 },
    @{ Row = 14; D = 0; H = 0;    I = "High"; J = "Urgent"; L = 1.11 },
    @{ Row = 15; D = 0; H = 0;    I = "High"; J = "Urgent"; L = 1.11 },
    @{ Row = 16; D = 0; H = 0;    I = "High"; J = "Urgent"; L = 1.19 },
    @{ Row = 17; D = 0; H = 0;    I = "High"; J = "Urgent"; L = 1.19 }
)

foreach ($r in $rows) {
    $wsForecast.Cells.Item($r.Row, 4).Value  = $r.D   # D: MyForecast
    $wsForecast.Cells.Item($r.Row, 8).Value  = $r.H   # H: Inventory Coverage
    $wsForecast.Cells.Item($r.Row, 9).Value  = $r.I   # I: Stockout Risk
    $wsForecast.Cells.Item($r.Row, 10).Value = $r.J   # J: Reorder Urgency
    $wsForecast.Cells.Item($r.Row, 12).Value = $r.L   # L: Seasonality Index
}

# Summary sheet roll-up figures, recomputed off the new MyForecast column.
# These cells hold text (not numbers) in the workbook, so a leading
# apostrophe is used to stop Excel from auto-coercing the assigned
# string into a numeric value.
$wsSummary.Range("B9").Value  = "'10"  # Total Forecast (16 Weeks)
$wsSummary.Range("B10").Value = "'7"   # Total Forecast (8 Weeks)
$wsSummary.Range("B11").Value = "'5"   # Total Forecast (4 Weeks)
$wsSummary.Range("B12").Value = "'3"   # Max Forecast
$wsSummary.Range("B14").Value = "'0"   # Min Forecast
